# Update gh-pages output values (F column "弹幕/热度" counts) for two
# worksheets: 展览 and 全部类型, rows 4, 5, 6, 9/11, 12/14.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 587
$ws1.Range("F5").Value = 533
$ws1.Range("F6").Value = 300
$ws1.Range("F9").Value = 7486
$ws1.Range("F12").Value = 29

# --- Sheet "全部类型" ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 587
$ws2.Range("F5").Value = 533
$ws2.Range("F6").Value = 300
$ws2.Range("F11").Value = 7486
$ws2.Range("F14").Value = 29
